$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.842714
$ws.Range("H2").Value = 41.685428
$ws.Range("I2").Value = 0.06176264451223276
$ws.Range("J2").Value = 0.04208443214243528
$ws.Range("M2").Value = 7.318981333333333
$ws.Range("N2").Value = 21.956944
$ws.Range("O2").Value = 0.1531761772116645
$ws.Range("P2").Value = 0.1531761772116645
$ws.Range("Q2").Value = 152.5474347020053
$ws.Range("R2").Value = 915.284608212032
$ws.Range("S2").Value = 0.009460565780866802
$ws.Range("T2").Value = 0.006446332435701936

$ws.Range("G3").Value = 20.842714
$ws.Range("H3").Value = 41.685428
$ws.Range("I3").Value = 0.06176264451223276
$ws.Range("J3").Value = 0.04208443214243528
$ws.Range("M3").Value = 32.599203
$ws.Range("N3").Value = 97.79760900000001
$ws.Range("O3").Value = 0.6822563234237459
$ws.Range("P3").Value = 0.6822563234237459
$ws.Range("Q3").Value = 679.4558647569421
$ws.Range("R3").Value = 4076.735188541652
$ws.Range("S3").Value = 0.04213795476984371
$ws.Range("T3").Value = 0.02871236994687401

$ws.Range("G4").Value = 20.842714
$ws.Range("H4").Value = 41.685428
$ws.Range("I4").Value = 0.06176264451223276
$ws.Range("J4").Value = 0.04208443214243528
$ws.Range("M4").Value = 7.863275333333333
$ws.Range("N4").Value = 23.589826
$ws.Range("O4").Value = 0.1645674993645896
$ws.Range("P4").Value = 0.1645674993645896
$ws.Range("Q4").Value = 163.8919988759213
$ws.Range("R4").Value = 983.351993255528
$ws.Range("S4").Value = 0.01016412396152224
$ws.Range("T4").Value = 0.006925729759859335

$ws.Range("I5").Value = 0.06445182531459281
$ws.Range("J5").Value = 0.0658752185158826
$ws.Range("M5").Value = 7.318981333333333
$ws.Range("N5").Value = 21.956944
$ws.Range("O5").Value = 0.1531761772116645
$ws.Range("P5").Value = 0.1531761772116645
$ws.Range("Q5").Value = 159.1894370982702
$ws.Range("R5").Value = 1432.704933884432
$ws.Range("S5").Value = 0.009872484216003312
$ws.Range("T5").Value = 0.01009051414524595

$ws.Range("I6").Value = 0.06445182531459281
$ws.Range("J6").Value = 0.0658752185158826
$ws.Range("M6").Value = 32.599203
$ws.Range("N6").Value = 97.79760900000001
$ws.Range("O6").Value = 0.6822563234237459
$ws.Range("P6").Value = 0.6822563234237459
$ws.Range("Q6").Value = 709.0397610098531
$ws.Range("R6").Value = 6381.357849088678
$ws.Range("S6").Value = 0.04397266537708361
$ws.Range("T6").Value = 0.04494378438938193

$ws.Range("I7").Value = 0.06445182531459281
$ws.Range("J7").Value = 0.0658752185158826
$ws.Range("M7").Value = 7.863275333333333
$ws.Range("N7").Value = 23.589826
$ws.Range("O7").Value = 0.1645674993645896
$ws.Range("P7").Value = 0.1645674993645896
$ws.Range("Q7").Value = 171.0279500729309
$ws.Range("R7").Value = 1539.251550656378
$ws.Range("S7").Value = 0.0106066757215059
$ws.Range("T7").Value = 0.01084091998125471

$ws.Range("G8").Value = 38.15032833333333
$ws.Range("H8").Value = 114.450985
$ws.Range("I8").Value = 0.1130498248393481
$ws.Range("J8").Value = 0.115546485737591
$ws.Range("M8").Value = 7.318981333333333
$ws.Range("N8").Value = 21.956944
$ws.Range("O8").Value = 0.1531761772116645
$ws.Range("P8").Value = 0.1531761772116645
$ws.Range("Q8").Value = 279.2215409322044
$ws.Range("R8").Value = 2512.99386838984
$ws.Range("S8").Value = 0.01731654000333961
$ws.Range("T8").Value = 0.01769896897552631

$ws.Range("G9").Value = 38.15032833333333
$ws.Range("H9").Value = 114.450985
$ws.Range("I9").Value = 0.1130498248393481
$ws.Range("J9").Value = 0.115546485737591
$ws.Range("M9").Value = 32.599203
$ws.Range("N9").Value = 97.79760900000001
$ws.Range("O9").Value = 0.6822563234237459
$ws.Range("P9").Value = 0.6822563234237459
$ws.Range("Q9").Value = 1243.670297854985
$ws.Range("R9").Value = 11193.03268069487
$ws.Range("S9").Value = 0.07712895785859207
$ws.Range("T9").Value = 0.07883232054386316

$ws.Range("G10").Value = 38.15032833333333
$ws.Range("H10").Value = 114.450985
$ws.Range("I10").Value = 0.1130498248393481
$ws.Range("J10").Value = 0.115546485737591
$ws.Range("M10").Value = 7.863275333333333
$ws.Range("N10").Value = 23.589826
$ws.Range("O10").Value = 0.1645674993645896
$ws.Range("P10").Value = 0.1645674993645896
$ws.Range("Q10").Value = 299.9865357420678
$ws.Range("R10").Value = 2699.87882167861
$ws.Range("S10").Value = 0.01860432697741638
$ws.Range("T10").Value = 0.01901519621820158

$ws.Range("G11").Value = 1.0325075
$ws.Range("H11").Value = 2.065015
$ws.Range("I11").Value = 0.003059601243807028
$ws.Range("J11").Value = 0.00208478088891425
$ws.Range("M11").Value = 7.318981333333333
$ws.Range("N11").Value = 21.956944
$ws.Range("O11").Value = 0.1531761772116645
$ws.Range("P11").Value = 0.1531761772116645
$ws.Range("Q11").Value = 7.556903119026666
$ws.Range("R11").Value = 45.34141871416
$ws.Range("S11").Value = 0.0004686580223184145
$ws.Range("T11").Value = 0.0003193387668878206

$ws.Range("G12").Value = 1.0325075
$ws.Range("H12").Value = 2.065015
$ws.Range("I12").Value = 0.003059601243807028
$ws.Range("J12").Value = 0.00208478088891425
$ws.Range("M12").Value = 32.599203
$ws.Range("N12").Value = 97.79760900000001
$ws.Range("O12").Value = 0.6822563234237459
$ws.Range("P12").Value = 0.6822563234237459
$ws.Range("Q12").Value = 33.6589215915225
$ws.Range("R12").Value = 201.953529549135
$ws.Range("S12").Value = 0.002087432295742503
$ws.Range("T12").Value = 0.001422354944414725

$ws.Range("G13").Value = 1.0325075
$ws.Range("H13").Value = 2.065015
$ws.Range("I13").Value = 0.003059601243807028
$ws.Range("J13").Value = 0.00208478088891425
$ws.Range("M13").Value = 7.863275333333333
$ws.Range("N13").Value = 23.589826
$ws.Range("O13").Value = 0.1645674993645896
$ws.Range("P13").Value = 0.1645674993645896
$ws.Range("Q13").Value = 8.118890756231666
$ws.Range("R13").Value = 48.71334453738999
$ws.Range("S13").Value = 0.0005035109257461108
$ws.Range("T13").Value = 0.0003430871776117045

$ws.Range("G14").Value = 239.6229553333334
$ws.Range("H14").Value = 718.868866
$ws.Range("I14").Value = 0.710068151739898
$ws.Range("J14").Value = 0.7257497274703861
$ws.Range("M14").Value = 7.318981333333333
$ws.Range("N14").Value = 21.956944
$ws.Range("O14").Value = 0.1531761772116645
$ws.Range("P14").Value = 0.1531761772116645
$ws.Range("Q14").Value = 1753.795937122834
$ws.Range("R14").Value = 15784.16343410551
$ws.Range("S14").Value = 0.1087655250432697
$ws.Range("T14").Value = 0.1111675688663211

$ws.Range("G15").Value = 239.6229553333334
$ws.Range("H15").Value = 718.868866
$ws.Range("I15").Value = 0.710068151739898
$ws.Range("J15").Value = 0.7257497274703861
$ws.Range("M15").Value = 32.599203
$ws.Range("N15").Value = 97.79760900000001
$ws.Range("O15").Value = 0.6822563234237459
$ws.Range("P15").Value = 0.6822563234237459
$ws.Range("Q15").Value = 7811.517364371267
$ws.Range("R15").Value = 70303.6562793414
$ws.Range("S15").Value = 0.4844484865863573
$ws.Range("T15").Value = 0.4951473407897312

$ws.Range("G16").Value = 239.6229553333334
$ws.Range("H16").Value = 718.868866
$ws.Range("I16").Value = 0.710068151739898
$ws.Range("J16").Value = 0.7257497274703861
$ws.Range("M16").Value = 7.863275333333333
$ws.Range("N16").Value = 23.589826
$ws.Range("O16").Value = 0.1645674993645896
$ws.Range("P16").Value = 0.1645674993645896
$ws.Range("Q16").Value = 1884.221273973035
$ws.Range("R16").Value = 16957.99146575732
$ws.Range("S16").Value = 0.116854140110271
$ws.Range("T16").Value = 0.1194348178143339

$ws.Range("G17").Value = 16.06600466666667
$ws.Range("H17").Value = 48.198014
$ws.Range("I17").Value = 0.04760795235012129
$ws.Range("J17").Value = 0.04865935524479072
$ws.Range("M17").Value = 7.318981333333333
$ws.Range("N17").Value = 21.956944
$ws.Range("O17").Value = 0.1531761772116645
$ws.Range("P17").Value = 0.1531761772116645
$ws.Range("Q17").Value = 117.5867882565796
$ws.Range("R17").Value = 1058.281094309216
$ws.Range("S17").Value = 0.007292404145866657
$ws.Range("T17").Value = 0.007453454021981398

$ws.Range("G18").Value = 16.06600466666667
$ws.Range("H18").Value = 48.198014
$ws.Range("I18").Value = 0.04760795235012129
$ws.Range("J18").Value = 0.04865935524479072
$ws.Range("M18").Value = 32.599203
$ws.Range("N18").Value = 97.79760900000001
$ws.Range("O18").Value = 0.6822563234237459
$ws.Range("P18").Value = 0.6822563234237459
$ws.Range("Q18").Value = 523.7389475276141
$ws.Range("R18").Value = 4713.650527748527
$ws.Range("S18").Value = 0.03248082653612663
$ws.Range("T18").Value = 0.03319815280948088

$ws.Range("G19").Value = 16.06600466666667
$ws.Range("H19").Value = 48.198014
$ws.Range("I19").Value = 0.04760795235012129
$ws.Range("J19").Value = 0.04865935524479072
$ws.Range("M19").Value = 7.863275333333333
$ws.Range("N19").Value = 23.589826
$ws.Range("O19").Value = 0.1645674993645896
$ws.Range("P19").Value = 0.1645674993645896
$ws.Range("Q19").Value = 126.3314182006182
$ws.Range("R19").Value = 1136.982763805564
$ws.Range("S19").Value = 0.007834721668127998
$ws.Range("T19").Value = 0.008007748413328438
